# Updated submission deadlines and instructions
#
# The supervisor-approval paragraph reads:
#   "... with your submission by May 13, 2016."
# where "May 13" (underlined) is split across two runs ("May 1" + "3").
# Update the deadline to "September 11, 2016", keeping "September 11"
# underlined as a single run, and re-anchor Word's auto-managed
# "_GoBack" bookmark to the newly edited span (its position after the
# edit follows the last place text was typed/replaced).

$d = $word.ActiveDocument
$newDate = "September 11"

# Locate the "May 13" span (it spans the "May 1" and "3" runs).
$dateRange = $d.Content
$found = $dateRange.Find.Execute("May 13", $true, $false, $false, $false,
                                  $false, $true, 1, $false, "", 0)
$start = $dateRange.Start

# Replace the old date text with the new one.
$dateRange.Text = $newDate

# Re-apply the underline so the new run matches the original formatting
# of the replaced "May 13" text.
$newRange = $d.Range($start, $start + $newDate.Length)
$newRange.Font.Underline = 1

# Re-anchor the "_GoBack" bookmark immediately after the newly typed text.
$markRange = $d.Range($start + $newDate.Length, $start + $newDate.Length)
[void]$d.Bookmarks.Add("_GoBack", $markRange)
